$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 13:22"

# Row 14 (Brasil)
$ws.Range("B14").Value = 46348
$ws.Range("C14").Value = 591
$ws.Range("E14").Value = 18096
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = 2934

# Row 79 (Republica de Macedonia)
$ws.Range("B79").Value = 1300
$ws.Range("C79").Value = 41
$ws.Range("D79").Value = 301
$ws.Range("E79").Value = 943
$ws.Range("F79").Value = 14

# Row 95 (Libano)
$ws.Range("B95").Value = 688
$ws.Range("C95").Value = 6
$ws.Range("D95").Value = 140
$ws.Range("E95").Value = 526

# Row 112 (Consejo Danes para los Refugiados)
$ws.Range("B112").Value = 377
$ws.Range("C112").Value = 18
$ws.Range("D112").Value = 47
$ws.Range("E112").Value = 305
